# add date select function.
# Updates the "activate" sheet with corrected/extra time entries, the
# "information" sheet's quality duration totals, and the "harvest" sheet's
# extra label rows for 2021-03-23.

$wb = $excel.ActiveWorkbook

$activate = $wb.Worksheets.Item("activate")
$information = $wb.Worksheets.Item("information")
$harvest = $wb.Worksheets.Item("harvest")

# ---------------------------------------------------------------------
# Sheet "activate": fix row 5 (endTime/duration) and append rows 6-11.
# ---------------------------------------------------------------------

$activate.Cells.Item(5, 3).Value = 44278.60706018518
$activate.Cells.Item(5, 4).Value = 348

$activateRows = @(
    @("2021-03-23", 44278.60728009259, 44278.63753472222, 2614, "paper"),
    @("2021-03-23", 44278.64493055556, 44278.68125,        3138, "paper"),
    @("2021-03-23", 44278.69215277778, 44278.7259837963,   2923, "paper"),
    @("2021-03-23", 44278.85123842592, 44278.8942824074,   3719, "paper"),
    @("2021-03-23", 44278.90400462963, 44278.92106481481,  1474, "paper"),
    @("2021-03-23", 44278.92135416667, 44278.93908564815,  1532, "think")
)

$r = 6
foreach ($row in $activateRows) {
    $activate.Cells.Item($r, 1).NumberFormat = "@"
    $activate.Cells.Item($r, 1).Value = $row[0]
    $activate.Cells.Item($r, 1).Style = "Normal"

    $activate.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $activate.Cells.Item($r, 2).Value = $row[1]

    $activate.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $activate.Cells.Item($r, 3).Value = $row[2]

    $activate.Cells.Item($r, 4).Value = $row[3]
    $activate.Cells.Item($r, 5).Value = $row[4]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet "information": update the duration totals for high/low quality.
# ---------------------------------------------------------------------

$information.Cells.Item(2, 3).Value = 23
$information.Cells.Item(4, 3).Value = 70

# ---------------------------------------------------------------------
# Sheet "harvest": append the day's label rows.
# ---------------------------------------------------------------------

$harvestRows = @(
    @("2021-03-23", "software", "convenience"),
    @("2021-03-23", "software", "research"),
    @("2021-03-23", "paper",    "note")
)

$r = 2
foreach ($row in $harvestRows) {
    $harvest.Cells.Item($r, 1).NumberFormat = "@"
    $harvest.Cells.Item($r, 1).Value = $row[0]
    $harvest.Cells.Item($r, 1).Style = "Normal"

    $harvest.Cells.Item($r, 2).Value = $row[1]
    $harvest.Cells.Item($r, 3).Value = $row[2]

    $r = $r + 1
}
